# Timesheet update: extend clinical prodrome data entry for row 16
# (week 19/06/2017 - 25/06/2017), per commit "finish first pass update of
# clinical prodrome papers".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Extend the "Work carried out" note for row 16 with a second line.
$ws.Range("D16").Value = "Data extracted from new prodrome pdfs" + [char]10 + "Extend clinical prodrome data"

# 2. The note now wraps onto two lines, so the row needs to be taller.
$ws.Rows.Item(16).RowHeight = 38

# 3. Hours worked that week increased from 5 to 13 (totals/formulas
#    downstream - E18, E20, E22 - recalculate automatically).
$ws.Range("E16").Value = 13
